$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L4").Value = 3471.45
$ws1.Range("N4").Value = 2116.73
$ws1.Range("D34").Value = 1286.22
$ws1.Range("D53").Value = "5 de 51"
$ws1.Range("N53").Value = "3 de 51"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F4").Value = 6976.72
$ws2.Range("F34").Value = 1286.22
$ws2.Range("F53").Value = 20391.41

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 6289.76
$ws3.Range("E3").Value = 21167.2476
$ws3.Range("F3").Value = 0.2290766747648058

$ws3.Range("D16").Value = 5910.1
$ws3.Range("E16").Value = 26831.35
$ws3.Range("F16").Value = 0.1805081937421831

$ws3.Range("D18").Value = 5352.79
$ws3.Range("E18").Value = -2152.79
$ws3.Range("F18").Value = 1.672746875

$ws3.Range("D19").Value = 21677.63
$ws3.Range("E19").Value = 72769.81064517914
$ws3.Range("F19").Value = 0.2295205656386041

$ws3.Columns.Item(5).ColumnWidth = 21.166666666666668
